$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 6 ("health drops below 1" / "hitting enemies until the health
# score drops below 1" / game-over description) is removed entirely; rows
# below it shift up by one.
$ws.Rows("6:6").Delete()

# The two previously-blank rows that are now rows 14 and 15 (after the
# shift) get new "health drops to 1" / "health drops to 0" test cases, and
# the row that is now row 16 gets a new "difficulty" edge-case test. Values
# are entered column-by-column (A's then B's then C/D's) to reproduce the
# shared-string insertion order of the authored workbook.
$ws.Range("A14").Value2 = "health drops to 1"
$ws.Range("A15").Value2 = "health drops to 0"

$ws.Range("B14").Value2 = "hitting enemies until the health score drops to 0"
$ws.Range("B15").Value2 = "hitting enemies until the health score drops to 0"

$ws.Range("C14").Value2 = "game continues to function as usual with all functions active, if hit one more time then health drops to 0"
$ws.Range("D14").Value2 = "game continues to function as usual with all functions active, if hit one more time then health drops to 0"

$ws.Range("C15").Value2 = "all functions cease and a game over screen is placed over the canvas with a game over message presenting the score and name of the player"
$ws.Range("D15").Value2 = "all functions cease and a game over screen is placed over the canvas with a game over message presenting the score and name of the player"

$ws.Range("A16").Value2 = "difficulty"
$ws.Range("B16").Value2 = "anything that is not ""easy"", ""medium"" or ""hard"""
$ws.Range("C16").Value2 = "unaccepted and reprompt"
$ws.Range("D16").Value2 = "unaccepted and reprompt"

# Restore the taller row heights those rows need to show their longer
# wrapped text (matches the rest of the sheet's hand-set heights).
$ws.Rows("14:14").RowHeight = 49.5
$ws.Rows("15:15").RowHeight = 66
$ws.Rows("16:16").RowHeight = 33

# Leave the same cell selected as in the saved workbook.
$ws.Range("D17").Select()
